# Generate Report for Handoff
# Adds two new localization-status rows (80badf6e-... and d237c79e-...)
# to the Overview sheet and each locale sheet (zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------
# Overview sheet: two new rows (6 and 7)
# ---------------------------------------------------------------------
$ws1.Range("A6").Value = "80badf6e-b6ac-4fc9-bab8-739414b53f04.md"
$ws1.Range("B6").Value = "Ready for handoff"
$ws1.Range("C6").Value = "Ready for handoff"
$ws1.Range("D6").Value = "2016-48-13 22:48:37"

$ws1.Range("A7").Value = "d237c79e-9d35-423d-8bce-5b95f1d6cbe1.md"
$ws1.Range("B7").Value = "Ready for handoff"
$ws1.Range("C7").Value = "Ready for handoff"
$ws1.Range("D7").Value = "2016-48-13 22:48:37"

$ws1.Hyperlinks.Add($ws1.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/80badf6ee6ac4fc9bab8739414b53f04badf6e0/e2e/80badf6e-b6ac-4fc9-bab8-739414b53f04.md", "", "", "80badf6e-b6ac-4fc9-bab8-739414b53f04.md")
$ws1.Hyperlinks.Add($ws1.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/d237c79e9d35423d8bce5b95f1d6cbe1d237c79/e2e/d237c79e-9d35-423d-8bce-5b95f1d6cbe1.md", "", "", "d237c79e-9d35-423d-8bce-5b95f1d6cbe1.md")

# ---------------------------------------------------------------------
# zh-cn sheet: two new rows (6 and 7)
# ---------------------------------------------------------------------
$ws2.Range("A6").Value = "80badf6e-b6ac-4fc9-bab8-739414b53f04.md"
$ws2.Range("B6").Value = ".md"
$ws2.Range("C6").Value = "Ready for handoff"
$ws2.Range("D6").Value = "80badf6e-b6ac-4fc9-bab8-739414b53f04.e93cacc08038aa03537ceeff43f86b0d0402e3ba.zh-cn.xlf"
$ws2.Range("E6").Value = "2016-03-13 22:48:33"
$ws2.Range("H6").Value = "0001-01-01 00:00:00"
$ws2.Range("I6").Value = "Include"

$ws2.Range("A7").Value = "d237c79e-9d35-423d-8bce-5b95f1d6cbe1.md"
$ws2.Range("B7").Value = ".md"
$ws2.Range("C7").Value = "Ready for handoff"
$ws2.Range("D7").Value = "d237c79e-9d35-423d-8bce-5b95f1d6cbe1.f41ff9e99bc6206078bd408f2e94f50f90306e32.zh-cn.xlf"
$ws2.Range("E7").Value = "2016-03-13 22:48:33"
$ws2.Range("H7").Value = "0001-01-01 00:00:00"
$ws2.Range("I7").Value = "Include"

$ws2.Hyperlinks.Add($ws2.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/80badf6ee6ac4fc9bab8739414b53f04badf6e0/e2e/80badf6e-b6ac-4fc9-bab8-739414b53f04.md", "", "", "80badf6e-b6ac-4fc9-bab8-739414b53f04.md")
$ws2.Hyperlinks.Add($ws2.Range("B6"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/80badf6ee6ac4fc9bab8739414b53f04badf6e0/e2e/80badf6e-b6ac-4fc9-bab8-739414b53f04.md", "", "", ".md")
$ws2.Hyperlinks.Add($ws2.Range("D6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e93cacc08038aa03537ceeff43f86b0d0402e3ba/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/80badf6e-b6ac-4fc9-bab8-739414b53f04.e93cacc08038aa03537ceeff43f86b0d0402e3ba.zh-cn.xlf", "", "", "80badf6e-b6ac-4fc9-bab8-739414b53f04.e93cacc08038aa03537ceeff43f86b0d0402e3ba.zh-cn.xlf")

$ws2.Hyperlinks.Add($ws2.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/d237c79e9d35423d8bce5b95f1d6cbe1d237c79/e2e/d237c79e-9d35-423d-8bce-5b95f1d6cbe1.md", "", "", "d237c79e-9d35-423d-8bce-5b95f1d6cbe1.md")
$ws2.Hyperlinks.Add($ws2.Range("B7"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/d237c79e9d35423d8bce5b95f1d6cbe1d237c79/e2e/d237c79e-9d35-423d-8bce-5b95f1d6cbe1.md", "", "", ".md")
$ws2.Hyperlinks.Add($ws2.Range("D7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f41ff9e99bc6206078bd408f2e94f50f90306e32/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/d237c79e-9d35-423d-8bce-5b95f1d6cbe1.f41ff9e99bc6206078bd408f2e94f50f90306e32.zh-cn.xlf", "", "", "d237c79e-9d35-423d-8bce-5b95f1d6cbe1.f41ff9e99bc6206078bd408f2e94f50f90306e32.zh-cn.xlf")

# ---------------------------------------------------------------------
# de-de sheet: two new rows (6 and 7)
# ---------------------------------------------------------------------
$ws3.Range("A6").Value = "80badf6e-b6ac-4fc9-bab8-739414b53f04.md"
$ws3.Range("B6").Value = ".md"
$ws3.Range("C6").Value = "Ready for handoff"
$ws3.Range("D6").Value = "80badf6e-b6ac-4fc9-bab8-739414b53f04.e93cacc08038aa03537ceeff43f86b0d0402e3ba.de-de.xlf"
$ws3.Range("E6").Value = "2016-03-13 22:48:37"
$ws3.Range("H6").Value = "0001-01-01 00:00:00"
$ws3.Range("I6").Value = "Include"

$ws3.Range("A7").Value = "d237c79e-9d35-423d-8bce-5b95f1d6cbe1.md"
$ws3.Range("B7").Value = ".md"
$ws3.Range("C7").Value = "Ready for handoff"
$ws3.Range("D7").Value = "d237c79e-9d35-423d-8bce-5b95f1d6cbe1.f41ff9e99bc6206078bd408f2e94f50f90306e32.de-de.xlf"
$ws3.Range("E7").Value = "2016-03-13 22:48:37"
$ws3.Range("H7").Value = "0001-01-01 00:00:00"
$ws3.Range("I7").Value = "Include"

$ws3.Hyperlinks.Add($ws3.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/80badf6ee6ac4fc9bab8739414b53f04badf6e0/e2e/80badf6e-b6ac-4fc9-bab8-739414b53f04.md", "", "", "80badf6e-b6ac-4fc9-bab8-739414b53f04.md")
$ws3.Hyperlinks.Add($ws3.Range("B6"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/80badf6ee6ac4fc9bab8739414b53f04badf6e0/e2e/80badf6e-b6ac-4fc9-bab8-739414b53f04.md", "", "", ".md")
$ws3.Hyperlinks.Add($ws3.Range("D6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e93cacc08038aa03537ceeff43f86b0d0402e3ba/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/80badf6e-b6ac-4fc9-bab8-739414b53f04.e93cacc08038aa03537ceeff43f86b0d0402e3ba.de-de.xlf", "", "", "80badf6e-b6ac-4fc9-bab8-739414b53f04.e93cacc08038aa03537ceeff43f86b0d0402e3ba.de-de.xlf")

$ws3.Hyperlinks.Add($ws3.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/d237c79e9d35423d8bce5b95f1d6cbe1d237c79/e2e/d237c79e-9d35-423d-8bce-5b95f1d6cbe1.md", "", "", "d237c79e-9d35-423d-8bce-5b95f1d6cbe1.md")
$ws3.Hyperlinks.Add($ws3.Range("B7"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/d237c79e9d35423d8bce5b95f1d6cbe1d237c79/e2e/d237c79e-9d35-423d-8bce-5b95f1d6cbe1.md", "", "", ".md")
$ws3.Hyperlinks.Add($ws3.Range("D7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f41ff9e99bc6206078bd408f2e94f50f90306e32/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/d237c79e-9d35-423d-8bce-5b95f1d6cbe1.f41ff9e99bc6206078bd408f2e94f50f90306e32.de-de.xlf", "", "", "d237c79e-9d35-423d-8bce-5b95f1d6cbe1.f41ff9e99bc6206078bd408f2e94f50f90306e32.de-de.xlf")

Write-Host "Report generated for handoff."
